$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "%" header in column D
$ws.Range("D1").Value = "%"

# New data values for "task completed" (B) and updated "total tasks" (C)
$ws.Range("B2").Value = 12
$ws.Range("B3").Value = 9
$ws.Range("C3").Value = 26

# "%" column formulas: task completed / total tasks
$ws.Range("D2").Formula = "=SUM(B2/C2)"
$ws.Range("D3:D5").Formula = "=SUM(B3/C3)"

# Percentage number format (0.0%) for the new column
$ws.Range("D2:D5").NumberFormat = "0.0%"

# Page setup matching the printed layout
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# Selection moved to F8
$ws.Range("F8").Select()
